$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Megha"

$ws.Range("A2").Select()
